# Add data for 2021-12-22 (workbook tracks "through Dec 13" -> "through Dec 14")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab and update the column header label to reflect the new cutoff date
$ws.Name = "Through 2021-12-14"
$ws.Range("B1").Value = "December 2021 (through December 14)"

# Updated / newly populated monthly counts per neighborhood
$ws.Range("B3").Value = 6
$ws.Range("D6").Value = 18
$ws.Range("BJ6").Value = 4
$ws.Range("AL7").Value = 5
$ws.Range("AX7").Value = 4
$ws.Range("B11").Value = 1
$ws.Range("BJ11").Value = 3
$ws.Range("AL23").Value = 2
$ws.Range("N24").Value = 4
$ws.Range("B34").Value = 4
$ws.Range("BV37").Value = 1
$ws.Range("Z38").Value = 1
$ws.Range("B41").Value = 5
$ws.Range("B42").Value = 3
$ws.Range("N43").Value = 1
$ws.Range("B54").Value = 1
$ws.Range("BV56").Value = 1
$ws.Range("B67").Value = 1
$ws.Range("AL67").Value = 1
$ws.Range("BJ67").Value = 2
$ws.Range("BJ75").Value = 1
$ws.Range("N82").Value = 3
$ws.Range("BJ82").Value = 1
$ws.Range("AX84").Value = 1
